# Refresh the scraped cryptocurrency quotes (Price / Volume(1h) columns, plus two
# rows whose Coin/Link/Price got reordered) to match the latest GitHub Actions run.
#
# Column D ("Price") is stored as text (see original inlineStr cells), so any new
# value that LOOKS like a plain number is written with a leading apostrophe -
# exactly what typing it into a text-bound cell in Excel does - forcing it to stay
# text instead of silently becoming a numeric cell. The apostrophe itself is never
# part of the stored value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.472.05'
$ws.Range("E2").Value = '  +0.38%  '

$ws.Range("D3").Value = '1.848.25'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("D4").Value = '''0.9983'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '''240.32'
$ws.Range("E5").Value = '  +0.04%  '

$ws.Range("D6").Value = '''0.6292'
$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").Value = '''0.9998'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '''0.07486'
$ws.Range("E8").Value = '  -1.50%  '

$ws.Range("D9").Value = '''0.2909'
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").Value = '''24.60'
$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("D11").Value = '''0.07743'
$ws.Range("E11").Value = '  -0.06%  '

$ws.Range("D12").Value = '1.902.38'
$ws.Range("E12").Value = '  +2.95%  '

$ws.Range("D13").Value = '''5.013'
$ws.Range("E13").Value = '  -0.02%  '

$ws.Range("D14").Value = '''0.6809'
$ws.Range("E14").Value = '  +0.14%  '

$ws.Range("D15").Value = '''0.00001043'
$ws.Range("E15").Value = '  -0.32%  '

$ws.Range("D16").Value = '''82.16'
$ws.Range("E16").Value = '  -1.08%  '

$ws.Range("D17").Value = '''6.240'
$ws.Range("E17").Value = '  +2.10%  '

$ws.Range("D18").Value = '29.477.37'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").Value = '''229.45'
$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").Value = '''0.9996'
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").Value = '''7.541'
$ws.Range("E22").Value = '  +1.17%  '

$ws.Range("D23").Value = '''0.9996'
$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").Value = '''159.47'
$ws.Range("E24").Value = '  +0.31%  '

$ws.Range("D25").Value = '''8.529'
$ws.Range("E25").Value = '  +1.07%  '

$ws.Range("D26").Value = '''0.1368'
$ws.Range("E26").Value = '  -1.61%  '

$ws.Range("D27").Value = '''17.54'
$ws.Range("E27").Value = '  -0.78%  '

$ws.Range("D28").Value = '''0.06553'
$ws.Range("E28").Value = '  +16.58%  '

$ws.Range("E29").Value = '  -1.30%  '

$ws.Range("D30").Value = '''1.487'
$ws.Range("E30").Value = '  +1.08%  '

$ws.Range("D31").Value = '''4.101'
$ws.Range("E31").Value = '  -0.18%  '

$ws.Range("D32").Value = '''4.108'
$ws.Range("E32").Value = '  +1.64%  '

$ws.Range("D33").Value = '''1.837'
$ws.Range("E33").Value = '  +0.37%  '

$ws.Range("D34").Value = '''1.146'
$ws.Range("E34").Value = '  -0.89%  '

$ws.Range("D35").Value = '''0.6994'
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("D36").Value = '''2.579'

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.01863'
$ws.Range("E37").Value = '  +1.99%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.265.79'
$ws.Range("E38").Value = '  +2.37%  '

$ws.Range("E39").Value = '  +4.33%  '

$ws.Range("D40").Value = '''6.828'
$ws.Range("E40").Value = '  +6.49%  '

$ws.Range("D41").Value = '''0.9339'
$ws.Range("E41").Value = '  +3.58%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").Value = '2.010.41'
$ws.Range("E43").Value = '  +0.22%  '

$ws.Range("D44").Value = '''101.36'
$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("D46").Value = '''1.735'
$ws.Range("E46").Value = '  +3.19%  '

$ws.Range("E47").Value = '  -0.79%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '''0.1164'
$ws.Range("E48").Value = '  +1.21%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000116'
$ws.Range("E49").Value = '  +1.23%  '

$ws.Range("D50").Value = '''8.997'
$ws.Range("E50").Value = '  +0.07%  '

$ws.Range("D51").Value = '''0.3959'
$ws.Range("E51").Value = '  -0.95%  '
